$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (nhl10@gmail.com) - this shifts rows 3,4,5 up to 2,3,4
$ws.Rows.Item(2).Delete()

# Clear the contents of the now-last row (A4/B4, previously row 5 = nhl13@gmail.com),
# but keep the cell's existing (hyperlink) style/formatting.
$ws.Hyperlinks.Item(2).Delete()
$ws.Range("A4:B4").ClearContents()

# Update the selection to match the target state
$ws.Range("F10").Select()
